# Junction_Flooding_201.xlsx edit:
#  1. Remove the last data row (row 6) - the sheet now ends at row 5.
#  2. Round row 5's hydraulic values down from 3 decimal places to 2
#     ("custom accuracy").
#  3. Narrow a subset of the data columns from width 8 to width 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete row 6 (also updates the sheet dimension to A1:AH5) ---
$ws.Rows(6).Delete()

# --- 2. Rewrite row 5's values with 2-decimal accuracy ---
$row5 = @{
    "B5"  = 21.84
    "C5"  = 16.49
    "D5"  = 0.61
    "E5"  = 47.88
    "F5"  = 39.31
    "G5"  = 16.97
    "H5"  = 64.7
    "I5"  = 26.6
    "J5"  = 12.27
    "K5"  = 17.83
    "L5"  = 19.28
    "M5"  = 20.42
    "N5"  = 5.53
    "O5"  = 17.27
    "P5"  = 24.49
    "Q5"  = 14.45
    "R5"  = 0.2
    "S5"  = 0.8100000000000001
    "T5"  = 255.92
    "U5"  = 48.3
    "V5"  = 15.94
    "W5"  = 32.43
    "X5"  = 17.38
    "Y5"  = 2.41
    "Z5"  = 32.18
    "AA5" = 14.02
    "AB5" = 12.76
    "AC5" = 15.06
    "AD5" = 20.48
    "AE5" = 0.33
    "AF5" = 58.61
    "AH5" = 19.79
}
foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# --- 3. Narrow the selected columns from width 8 to width 7 ---
# Excel's ColumnWidth property is offset from the raw OOXML <col width>
# attribute by 5/6 of a character (the default column-width padding), so
# subtract that offset to land exactly on the integer width stored in XML.
$narrowCols = @("B","C","G","J","K","L","M","O","Q","V","X","AA","AB","AC","AD","AH")
$targetWidth = 7 - (5 / 6)
foreach ($col in $narrowCols) {
    $ws.Columns($col).ColumnWidth = $targetWidth
}
